# New weekly price report: a fresh record for the latest week is inserted
# at row 64 (pushing the existing rows 64-102 down to 65-103), matching the
# "Fruta / hortaliza, semanal" update pattern used across these sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 64; everything from the old
# row 64 downward shifts down by one (old row 102 becomes row 103).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with this week's record.
$ws.Range("A64").Value = 8
$ws.Range("B64").Value = "Terminal La Palmera de La Serena"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44957
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = 100114007
$ws.Range("G64").Value = "Jengibre"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 400
$ws.Range("K64").Value = 23000
$ws.Range("L64").Value = 24000
$ws.Range("M64").Value = 23500
$ws.Range("N64").Value = "$/caja 13 kilos"
$ws.Range("O64").Value = "Perú"
$ws.Range("P64").Value = 1808
$ws.Range("Q64").Value = 13
$ws.Range("R64").Value = "Hortaliza"
